$wb = $excel.ActiveWorkbook

# --- "Texts" sheet: insert a new "Restart Game" text column between
#     "Score" (col W) and "Points" (col X, now shifted to Y) ---
$wsTexts = $wb.Worksheets.Item("Texts")
$wsTexts.Range("X1").EntireColumn.Insert()
$wsTexts.Range("X1").Value = "Restart Game"
$wsTexts.Range("X2").Value = "Reiniciar o jogo"
$wsTexts.Range("X3").Value = "Restart game"
$wsTexts.Range("X4").Value = "Reiniciar el juego"
$wsTexts.Range("X11").Select()

# --- "Times" sheet: rename the header labels ("pista"/"tempo" -> "lane"/"time") ---
$wsTimes = $wb.Worksheets.Item("Times")
$wsTimes.Range("A1").Value = "lane"
$wsTimes.Range("B1").Value = "time"
$wsTimes.Range("E11").Select()

# --- "Debug" sheet becomes the active tab (was "Export") ---
$wsDebug = $wb.Worksheets.Item("Debug")
$wsDebug.Activate()
$wsDebug.Range("B8").Select()
